# updated URL for Aseul
# The "Bye Bye Summer" / Aseul row (row 11) pointed at an old YouTube video
# (https://youtu.be/4Ti34K4OhCY). Update the displayed YouTubeURL cell and
# the corresponding EmbedCode <iframe> src to the new video
# (https://youtu.be/Vs-J3-1YQ-I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = "https://youtu.be/Vs-J3-1YQ-I"
$ws.Range("G11").Value = '<iframe width="560" height="315" src="https://www.youtube.com/embed/Vs-J3-1YQ-I" frameborder="0" allow="accelerometer; autoplay; clipboard-write; encrypted-media; gyroscope; picture-in-picture" allowfullscreen></iframe>'
